# Reorder the "kitchens_2" / "living_rooms_1" block columns (B and D) on the
# active sheet: header labels swap, and the one-hot "1" markers in rows 3, 4
# and 6 move along with their column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header labels: swap B1 and D1
$ws.Range("B1").Value = "living_rooms_1"
$ws.Range("D1").Value = "kitchens_2"

# Row 3: the "1" marker moves from column A to column D
$ws.Range("A3").Value = 0
$ws.Range("D3").Value = 1

# Row 4: the "1" marker moves from column D to column B
$ws.Range("B4").Value = 1
$ws.Range("D4").Value = 0

# Row 6: the "1" marker moves from column B to column A
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 0
